$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new value, whether the value should be forced as text
# (because it would otherwise be auto-recognized by Excel as a number).
$updates = @(
    @(2, "B", 'Bitcoin', $false),
    @(2, "C", 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', $false),
    @(2, "D", '91.185.18', $false),
    @(2, "E", '  +3.30%  ', $false),
    @(3, "B", 'Ethereum', $false),
    @(3, "C", 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', $false),
    @(3, "D", '3.134.60', $false),
    @(3, "E", '  +1.58%  ', $false),
    @(4, "B", 'TetherUSD', $false),
    @(4, "C", 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', $false),
    @(4, "D", '1.00', $true),
    @(4, "E", '  -0.02%  ', $false),
    @(5, "B", 'Solana', $false),
    @(5, "C", 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', $false),
    @(5, "D", '220.34', $true),
    @(5, "E", '  +4.97%  ', $false),
    @(6, "B", 'BNB', $false),
    @(6, "C", 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', $false),
    @(6, "D", '623.80', $true),
    @(6, "E", '  +0.06%  ', $false),
    @(7, "B", 'Dogecoin', $false),
    @(7, "C", 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', $false),
    @(7, "D", '0.378', $true),
    @(7, "E", '  +1.76%  ', $false),
    @(8, "B", 'XRP', $false),
    @(8, "C", 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', $false),
    @(8, "D", '0.914', $true),
    @(8, "E", '  +11.33%  ', $false),
    @(9, "B", 'USDC', $false),
    @(9, "C", 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', $false),
    @(9, "D", '1.00', $true),
    @(9, "E", '  -0.03%  ', $false),
    @(10, "B", 'LidoStakedEther', $false),
    @(10, "C", 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', $false),
    @(10, "D", '3.128.53', $false),
    @(10, "E", '  +1.42%  ', $false),
    @(11, "B", 'Cardano', $false),
    @(11, "C", 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', $false),
    @(11, "D", '0.746', $true),
    @(11, "E", '  +25.77%  ', $false),
    @(12, "B", 'TRON', $false),
    @(12, "C", 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', $false),
    @(12, "D", '0.190', $true),
    @(12, "E", '  +6.52%  ', $false),
    @(13, "B", 'ShibaInu', $false),
    @(13, "C", 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', $false),
    @(13, "D", '0.0000254', $true),
    @(13, "E", '  +6.65%  ', $false),
    @(14, "B", 'Avalanche', $false),
    @(14, "C", 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', $false),
    @(14, "D", '34.18', $true),
    @(14, "E", '  +8.05%  ', $false),
    @(15, "B", 'Toncoin', $false),
    @(15, "C", 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', $false),
    @(15, "D", '5.42', $true),
    @(15, "E", '  +2.59%  ', $false),
    @(16, "B", 'WrappedBTC', $false),
    @(16, "C", 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', $false),
    @(16, "D", '91.027.51', $false),
    @(16, "E", '  +3.47%  ', $false),
    @(17, "B", 'WrappedliquidstakedEther2.0', $false),
    @(17, "C", 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', $false),
    @(17, "D", '3.707.96', $false),
    @(17, "E", '  +1.48%  ', $false),
    @(18, "B", 'WrappedEther', $false),
    @(18, "C", 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', $false),
    @(18, "D", '3.153.35', $false),
    @(18, "E", '  +2.26%  ', $false),
    @(19, "B", 'SuiNetwork', $false),
    @(19, "C", 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', $false),
    @(19, "D", '3.85', $true),
    @(19, "E", '  +20.29%  ', $false),
    @(20, "B", 'PEPE', $false),
    @(20, "C", 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', $false),
    @(20, "D", '0.0000228', $true),
    @(20, "E", '  +7.50%  ', $false),
    @(21, "B", 'Chainlink', $false),
    @(21, "C", 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', $false),
    @(21, "D", '14.19', $true),
    @(21, "E", '  +8.40%  ', $false),
    @(22, "B", 'BitcoinCash', $false),
    @(22, "C", 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', $false),
    @(22, "D", '432.81', $true),
    @(22, "E", '  +2.76%  ', $false),
    @(23, "B", 'Uniswap', $false),
    @(23, "C", 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', $false),
    @(23, "D", '8.78', $true),
    @(23, "E", '  +7.45%  ', $false),
    @(24, "B", 'Polkadot', $false),
    @(24, "C", 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', $false),
    @(24, "D", '5.17', $true),
    @(24, "E", '  +6.62%  ', $false),
    @(25, "B", 'NEARProtocol', $false),
    @(25, "C", 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', $false),
    @(25, "D", '6.10', $true),
    @(25, "E", '  +12.16%  ', $false),
    @(26, "B", 'Aptos', $false),
    @(26, "C", 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', $false),
    @(26, "D", '12.43', $true),
    @(26, "E", '  +8.40%  ', $false),
    @(27, "B", 'Litecoin', $false),
    @(27, "C", 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', $false),
    @(27, "D", '83.84', $true),
    @(27, "E", '  +2.66%  ', $false),
    @(28, "B", 'Dai', $false),
    @(28, "C", 'https://coinranking.com/coin/MoTuySvg7+dai-dai', $false),
    @(28, "D", '0.999', $true),
    @(28, "E", '  -0.09%  ', $false),
    @(29, "B", 'Cronos', $false),
    @(29, "C", 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', $false),
    @(29, "D", '0.167', $true),
    @(29, "E", '  +8.48%  ', $false),
    @(30, "B", 'InternetComputer(DFINITY)', $false),
    @(30, "C", 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', $false),
    @(30, "D", '9.06', $true),
    @(30, "E", '  +12.76%  ', $false),
    @(31, "B", 'Bittensor', $false),
    @(31, "C", 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', $false),
    @(31, "D", '530.04', $true),
    @(31, "E", '  +5.08%  ', $false),
    @(32, "B", 'Binance-PegBSC-USD', $false),
    @(32, "C", 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', $false),
    @(32, "D", '0.912', $true),
    @(32, "E", '  -16.17%  ', $false),
    @(33, "B", 'dogwifhat', $false),
    @(33, "C", 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', $false),
    @(33, "D", '3.88', $true),
    @(33, "E", '  +9.43%  ', $false),
    @(34, "B", 'RenderToken', $false),
    @(34, "C", 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render', $false),
    @(34, "D", '7.39', $true),
    @(34, "E", '  +12.04%  ', $false),
    @(35, "B", 'Fetch.AI', $false),
    @(35, "C", 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', $false),
    @(35, "D", '1.31', $true),
    @(35, "E", '  +6.70%  ', $false),
    @(36, "B", 'Kaspa', $false),
    @(36, "C", 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', $false),
    @(36, "D", '0.141', $true),
    @(36, "E", '  +7.19%  ', $false),
    @(37, "B", 'EthereumClassic', $false),
    @(37, "C", 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', $false),
    @(37, "D", '23.47', $true),
    @(37, "E", '  +5.54%  ', $false),
    @(38, "B", 'PancakeSwap', $false),
    @(38, "C", 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', $false),
    @(38, "D", '1.86', $true),
    @(38, "E", '  +2.78%  ', $false),
    @(39, "B", 'WhiteBITCoin', $false),
    @(39, "C", 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', $false),
    @(39, "D", '22.30', $true),
    @(39, "E", '  +0.39%  ', $false),
    @(40, "B", 'FirstDigitalUSD', $false),
    @(40, "C", 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', $false),
    @(40, "D", '1.00', $true),
    @(40, "E", '  -0.03%  ', $false),
    @(41, "B", 'Stellar', $false),
    @(41, "C", 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', $false),
    @(41, "D", '0.144', $true),
    @(41, "E", '  +6.92%  ', $false),
    @(42, "B", 'USDe', $false),
    @(42, "C", 'https://coinranking.com/coin/exbfr2U-0+usde-usde', $false),
    @(42, "D", '1.00', $true),
    @(42, "E", '  +0.02%  ', $false),
    @(43, "B", 'Hedera', $false),
    @(43, "C", 'https://coinranking.com/coin/jad286TjB+hedera-hbar', $false),
    @(43, "D", '0.0771', $true),
    @(43, "E", '  +16.96%  ', $false),
    @(44, "B", 'PolygonEcosystemToken', $false),
    @(44, "C", 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', $false),
    @(44, "D", '0.378', $true),
    @(44, "E", '  +5.63%  ', $false),
    @(45, "B", 'Stacks', $false),
    @(45, "C", 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', $false),
    @(45, "D", '1.92', $true),
    @(45, "E", '  +5.79%  ', $false),
    @(46, "B", 'Monero', $false),
    @(46, "C", 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', $false),
    @(46, "D", '143.67', $true),
    @(46, "E", '  -3.80%  ', $false),
    @(47, "B", 'OKB', $false),
    @(47, "C", 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', $false),
    @(47, "D", '44.17', $true),
    @(47, "E", '  +1.61%  ', $false),
    @(48, "B", 'ImmutableX', $false),
    @(48, "C", 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', $false),
    @(48, "D", '1.30', $true),
    @(48, "E", '  +11.32%  ', $false),
    @(49, "B", 'FLOKI', $false),
    @(49, "C", 'https://coinranking.com/coin/fmHk13Rqw+floki-floki', $false),
    @(49, "D", '0.000266', $true),
    @(49, "E", '  +24.58%  ', $false),
    @(50, "B", 'Aave', $false),
    @(50, "C", 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', $false),
    @(50, "D", '168.59', $true),
    @(50, "E", '  +8.04%  ', $false),
    @(51, "B", 'Filecoin', $false),
    @(51, "C", 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', $false),
    @(51, "D", '4.18', $true),
    @(51, "E", '  +7.42%  ', $false)

)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $forceText = $u[3]
    $addr = "$col$row"
    $rng = $ws.Range($addr)
    if ($forceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = "Normal"
    } else {
        $rng.Value = $val
    }
}
